# Chỉnh sửa model Product, xóa bảng Price, Seed lại data
#
# The old sheet had two columns:
#   A: Id   (header "Id" + numeric ids 1..13)
#   B: Name (header "Name" + product name strings)
#
# The "Price" table concept (column A / the numeric Id column) is removed,
# and the remaining Name data is reseeded into column A by deleting the
# old column A outright (Excel shifts column B left into column A).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old Id column entirely; Name (old column B) shifts into column A.
$ws.Columns.Item(1).Delete()

# Match the new selection recorded for the sheet after the edit.
$ws.Range("C3").Select()
